# CodeSystem-CSSexoBiologico.xlsx — "version final sin errores"
#
# Semantic changes applied (per the authoritative OOXML diff):
#   1. Metadata sheet: "Version" row's value changes from "0.4.0" to "0.7.0".
#   2. Metadata sheet: the "Jurisdiction" / "Chile" row is removed entirely,
#      shifting every row below it up by one (dimension A1:B22 -> A1:B21).
#   3. Concepts sheet is untouched (its apparent shared-string index churn in
#      the diff is purely a side effect of the two strings above disappearing
#      from the shared-string table; the cell contents are identical).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Bump the Version value (row 3, column B) to 0.7.0
$ws.Range("B3").Value = "0.7.0"

# 2) Delete the whole "Jurisdiction" / "Chile" row (row 11)
$ws.Rows("11:11").Delete()
